$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2292
$ws1.Range("F7").Value = 323
$ws1.Range("F10").Value = 262
$ws1.Range("F11").Value = 476
$ws1.Range("C12").Value = "杭州·第十届次元鹿角动漫游戏展（取消）"
$ws1.Range("G12").Value = "不可售"
$ws1.Range("F16").Value = 8080
$ws1.Range("F19").Value = 226
$ws1.Range("F27").Value = 1868
$ws1.Range("F28").Value = 601
$ws1.Range("F30").Value = 1691
$ws1.Range("F31").Value = 238
$ws1.Range("F34").Value = 7
$ws1.Range("F39").Value = 195
$ws1.Range("F40").Value = 362
$ws1.Range("F42").Value = 231

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 2

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 2292
$ws4.Range("F9").Value = 323
$ws4.Range("F13").Value = 262
$ws4.Range("C15").Value = "杭州·第十届次元鹿角动漫游戏展（取消）"
$ws4.Range("G15").Value = "不可售"
$ws4.Range("F18").Value = 8080
$ws4.Range("F22").Value = 226
$ws4.Range("F30").Value = 1868
$ws4.Range("F31").Value = 602
$ws4.Range("F33").Value = 1691
$ws4.Range("F34").Value = 238
$ws4.Range("F37").Value = 7
$ws4.Range("F42").Value = 195
$ws4.Range("F43").Value = 362
$ws4.Range("F44").Value = 2
$ws4.Range("F49").Value = 231
